# Update countries & provincias Spain
# Applies the data refresh captured in the commit: the "Pais" (countries)
# sheet is resorted by total-cases as new numbers come in (a few countries
# leap-frog their neighbours in the table), several rows get refreshed
# totals, and the "last updated" timestamp moves forward.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-Row($row, $country, $total, $nuevos, $activos, $recuperados, $criticos, $muertesHoy, $muertes) {
    $ws.Cells.Item($row, 1).Value = $country
    $ws.Cells.Item($row, 2).Value = $total
    $ws.Cells.Item($row, 3).Value = $nuevos
    $ws.Cells.Item($row, 4).Value = $activos
    $ws.Cells.Item($row, 5).Value = $recuperados
    $ws.Cells.Item($row, 6).Value = $criticos
    $ws.Cells.Item($row, 7).Value = $muertesHoy
    $ws.Cells.Item($row, 8).Value = $muertes
}

# --- timestamp banner (A1) -------------------------------------------------
$ws.Range("A1").Value = "Datos actualizados a 30 de Abril de 2020 a las 02:22"

# --- pure numeric refreshes (no reordering) --------------------------------
Set-Row 4   "Estados Unidos"      1063854 28089 147114 855105 18671 2369 61635
Set-Row 9   "Alemania"            161539  1627  120400 34672  2415  153  6467
Set-Row 32  "Japon"               13895   159   2368   11114  306   19   413
Set-Row 117 "Reunion"             420     2     300    120    2     0    0
Set-Row 146 "Trinidad yTobago"    116     0     71     37     0     0    8
Set-Row 166 "Polinesia Francesa"  58      0     50     8      1     0    0

# --- Ecuador overtakes Portugal (rows 21/22 swap order + refresh) ---------
Set-Row 21 "Ecuador"  24675 417 1557 22235 146 12 883
Set-Row 22 "Portugal" 24505 183 1470 22062 169 25 973

# --- Niger overtakes Crucero (rows 101/102 swap order + refresh) ----------
Set-Row 101 "Niger"   713 4 435 246 0 1 32
Set-Row 102 "Crucero" 712 0 645 54  4 0 13

# --- Gabon jumps ahead of Vietnam & Paraguay (rows 129-131 rotate) --------
Set-Row 129 "Gabon"    276 38 67  206 1 0 3
Set-Row 130 "Vietnam"  270 0  222 48  8 0 0
Set-Row 131 "Paraguay" 239 9  102 128 1 0 9

# --- San Vicente y las Granadinas overtakes Namibia (rows 191/192 swap) ---
Set-Row 191 "San Vicente y las Granadinas" 16 1 8 8 0 0 0
Set-Row 192 "Namibia"                       16 0 8 8 0 0 0
